$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new first column (shifts existing A:D -> B:E)
$ws.Columns.Item(1).EntireColumn.Insert()

# New column A header/value - tab name labels
$ws.Range("A1").Value = "TabName"
$ws.Range("A2").Value = "CasesTab"

# Updated query text for the (now) B2 and C2 cells - new Cypher queries
$casesTabQuery = "MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)`n    WHERE c.race = `"NOT_REPORTED`"`nWITH DISTINCT c, a, ct`nRETURN `n    COALESCE(c.case_id, '') AS ``Case ID``,`n    COALESCE(ct.clinical_trial_designation, '') AS ``Trial Code``,`n    COALESCE(a.arm_id, '') AS ``Arm``,`n    COALESCE(a.arm_drug, '') AS ``Arm Treatment``,`n    COALESCE(c.disease, '') AS ``Diagnosis``,`n    COALESCE(c.gender, '') AS ``Gender``,`n    COALESCE(c.race, '') AS ``Race``,`n    COALESCE(c.ethnicity, '') AS ``Ethnicity``"

$statQuery = "MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)`n    WHERE c.race = `"NOT_REPORTED`"`nOPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)`nRETURN `n    COUNT(DISTINCT f) AS number_of_files,`n    COUNT(DISTINCT c.case_id) AS number_of_cases,`n    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials"

$ws.Range("B2").Value = $casesTabQuery
$ws.Range("C2").Value = $statQuery

# Row 2 grew taller to fit the longer wrapped text
$ws.Rows.Item(2).RowHeight = 174

# New column A sizing (best-fit width for short labels)
$ws.Columns.Item(1).ColumnWidth = 7.92

# Update selection to match final state
$ws.Range("B2").Select()
